$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This workbook lists species-observation records, one per row (row number
# is a stable key). The edit re-orders which observation records data sits
# in which row (several small permutations of row groups), while row 1
# (headers) and all unaffected rows stay untouched. We reproduce this by
# writing, cell by cell, the new contents that belong in each affected row.

# --- Row 2 ---
$ws.Range("A2").Value2 = 130789503
$ws.Range("B2").Value2 = 79243
$ws.Range("E2").Value2 = 6425
$ws.Range("F2").Value2 = "Garnlav"
$ws.Range("G2").Value2 = "Alectoria sarmentosa"
$ws.Range("H2").Value2 = "(Ach.) Ach."
$ws.Range("K2").ClearContents() | Out-Null
$ws.Range("L2").ClearContents() | Out-Null
$ws.Range("M2").ClearContents() | Out-Null
$ws.Range("N2").ClearContents() | Out-Null
$ws.Range("Q2").Value2 = 490685
$ws.Range("R2").Value2 = 6763486
$ws.Range("Z2").Value2 = "12:18"
$ws.Range("AB2").Value2 = "12:18"
$ws.Range("AF2").ClearContents() | Out-Null

# --- Row 3 ---
$ws.Range("A3").Value2 = 130789515
$ws.Range("B3").Value2 = 57881
$ws.Range("E3").Value2 = 100049
$ws.Range("F3").Value2 = "Spillkråka"
$ws.Range("G3").Value2 = "Dryocopus martius"
$ws.Range("H3").Value2 = "(Linnaeus, 1758)"
$ws.Range("K3").ClearContents() | Out-Null
$ws.Range("L3").ClearContents() | Out-Null
$ws.Range("M3").Value2 = "gammalt bo"
$ws.Range("N3").ClearContents() | Out-Null
$ws.Range("Q3").Value2 = 490723
$ws.Range("R3").Value2 = 6763501
$ws.Range("Z3").Value2 = "12:28"
$ws.Range("AB3").Value2 = "12:28"
$ws.Range("AF3").ClearContents() | Out-Null

# --- Row 7 ---
$ws.Range("A7").Value2 = 130789501
$ws.Range("B7").Value2 = 79243
$ws.Range("E7").Value2 = 6425
$ws.Range("F7").Value2 = "Garnlav"
$ws.Range("G7").Value2 = "Alectoria sarmentosa"
$ws.Range("H7").Value2 = "(Ach.) Ach."
$ws.Range("M7").ClearContents() | Out-Null
$ws.Range("P7").Value2 = "Kråkbackarna, Dlr"
$ws.Range("Q7").Value2 = 490713
$ws.Range("R7").Value2 = 6763507
$ws.Range("Z7").Value2 = "12:36"
$ws.Range("AB7").Value2 = "12:36"
$ws.Range("AF7").ClearContents() | Out-Null
$ws.Range("AW7").Value2 = "Bo karlstens"
$ws.Range("AX7").Value2 = "Bo karlstens, Håkan Thenander"

# --- Row 8 ---
$ws.Range("A8").Value2 = 130752192
$ws.Range("B8").Value2 = 57881
$ws.Range("E8").Value2 = 100049
$ws.Range("F8").Value2 = "Spillkråka"
$ws.Range("G8").Value2 = "Dryocopus martius"
$ws.Range("H8").Value2 = "(Linnaeus, 1758)"
$ws.Range("M8").Value2 = "äldre spår"
$ws.Range("P8").Value2 = "Truppan, Dlr"
$ws.Range("Q8").Value2 = 490715
$ws.Range("R8").Value2 = 6763290
$ws.Range("Z8").Value2 = "11:43"
$ws.Range("AB8").Value2 = "11:43"
$ws.Range("AF8").ClearContents() | Out-Null
$ws.Range("AW8").Value2 = "Håkan Thenander"
$ws.Range("AX8").Value2 = "Håkan Thenander, Bo karlstens"

# --- Row 9 ---
$ws.Range("A9").Value2 = 130789509
$ws.Range("Q9").Value2 = 490693
$ws.Range("R9").Value2 = 6763417
$ws.Range("Z9").Value2 = "12:05"
$ws.Range("AB9").Value2 = "12:05"

# --- Row 12 ---
$ws.Range("A12").Value2 = 130751852
$ws.Range("P12").Value2 = "Truppan, Dlr"
$ws.Range("Q12").Value2 = 490760
$ws.Range("R12").Value2 = 6763211
$ws.Range("Z12").Value2 = "11:43"
$ws.Range("AB12").Value2 = "11:43"
$ws.Range("AF12").ClearContents() | Out-Null
$ws.Range("AW12").Value2 = "Håkan Thenander"
$ws.Range("AX12").Value2 = "Håkan Thenander, Bo karlstens"

# --- Row 14 ---
$ws.Range("A14").Value2 = 130789504
$ws.Range("P14").Value2 = "Kråkbackarna, Dlr"
$ws.Range("Q14").Value2 = 490686
$ws.Range("R14").Value2 = 6763480
$ws.Range("Z14").Value2 = "12:16"
$ws.Range("AB14").Value2 = "12:16"
$ws.Range("AF14").ClearContents() | Out-Null
$ws.Range("AW14").Value2 = "Bo karlstens"
$ws.Range("AX14").Value2 = "Bo karlstens, Håkan Thenander"

# --- Row 18 ---
$ws.Range("A18").Value2 = 130789512
$ws.Range("Q18").Value2 = 490686
$ws.Range("R18").Value2 = 6763382
$ws.Range("Z18").Value2 = "12:01"
$ws.Range("AB18").Value2 = "12:01"
$ws.Range("AF18").ClearContents() | Out-Null
$ws.Range("AW18").Value2 = "Bo karlstens"
$ws.Range("AX18").Value2 = "Bo karlstens, Håkan Thenander"

# --- Row 19 ---
$ws.Range("A19").Value2 = 130789514
$ws.Range("Q19").Value2 = 490800
$ws.Range("R19").Value2 = 6763195
$ws.Range("Z19").Value2 = "11:44"
$ws.Range("AB19").Value2 = "11:44"

# --- Row 20 ---
$ws.Range("A20").Value2 = 130752453
$ws.Range("Q20").Value2 = 490682
$ws.Range("R20").Value2 = 6763392
$ws.Range("Z20").Value2 = "11:43"
$ws.Range("AB20").Value2 = "11:43"
$ws.Range("AF20").ClearContents() | Out-Null
$ws.Range("AW20").Value2 = "Håkan Thenander"
$ws.Range("AX20").Value2 = "Håkan Thenander, Bo karlstens"

# --- Row 24 ---
$ws.Range("A24").Value2 = 130752740
$ws.Range("P24").Value2 = "Kråkbackarna, Dlr"
$ws.Range("Q24").Value2 = 490682
$ws.Range("R24").Value2 = 6763461

# --- Row 25 ---
$ws.Range("A25").Value2 = 130789507
$ws.Range("Q25").Value2 = 490706
$ws.Range("R25").Value2 = 6763438
$ws.Range("Z25").Value2 = "12:09"
$ws.Range("AB25").Value2 = "12:09"
$ws.Range("AF25").ClearContents() | Out-Null
$ws.Range("AW25").Value2 = "Bo karlstens"
$ws.Range("AX25").Value2 = "Bo karlstens, Håkan Thenander"

# --- Row 26 ---
$ws.Range("A26").Value2 = 130752001
$ws.Range("P26").Value2 = "Truppan, Dlr"
$ws.Range("Q26").Value2 = 490746
$ws.Range("R26").Value2 = 6763219
$ws.Range("Z26").Value2 = "11:43"
$ws.Range("AB26").Value2 = "11:43"
$ws.Range("AF26").ClearContents() | Out-Null
$ws.Range("AW26").Value2 = "Håkan Thenander"
$ws.Range("AX26").Value2 = "Håkan Thenander, Bo karlstens"

# --- Row 27 ---
$ws.Range("A27").Value2 = 130751938
$ws.Range("P27").Value2 = "Truppan, Dlr"
$ws.Range("Q27").Value2 = 490749
$ws.Range("R27").Value2 = 6763201

# --- Row 28 ---
$ws.Range("A28").Value2 = 130752569
$ws.Range("P28").Value2 = "Kråkbackarna, Dlr"
$ws.Range("Q28").Value2 = 490661
$ws.Range("R28").Value2 = 6763445

# --- Row 32 ---
$ws.Range("A32").Value2 = 130752842
$ws.Range("Q32").Value2 = 490660
$ws.Range("R32").Value2 = 6763462

# --- Row 33 ---
$ws.Range("A33").Value2 = 130752874
$ws.Range("Q33").Value2 = 490655
$ws.Range("R33").Value2 = 6763493

# --- Row 34 ---
$ws.Range("A34").Value2 = 130753455
$ws.Range("B34").Value2 = 57884
$ws.Range("E34").Value2 = 100109
$ws.Range("F34").Value2 = "Tretåig hackspett"
$ws.Range("G34").Value2 = "Picoides tridactylus"
$ws.Range("H34").Value2 = "(Linnaeus, 1758)"
$ws.Range("M34").Value2 = "färska spår"
$ws.Range("Q34").Value2 = 490770
$ws.Range("R34").Value2 = 6763512
$ws.Range("AC34").Value2 = "Bild 4 till 6"

# --- Row 35 ---
$ws.Range("A35").Value2 = 130753055
$ws.Range("B35").Value2 = 79243
$ws.Range("E35").Value2 = 6425
$ws.Range("F35").Value2 = "Garnlav"
$ws.Range("G35").Value2 = "Alectoria sarmentosa"
$ws.Range("H35").Value2 = "(Ach.) Ach."
$ws.Range("M35").ClearContents() | Out-Null
$ws.Range("Q35").Value2 = 490658
$ws.Range("R35").Value2 = 6763532
$ws.Range("AC35").ClearContents() | Out-Null

# --- Row 36 ---
$ws.Range("A36").Value2 = 130789513
$ws.Range("Q36").Value2 = 490777
$ws.Range("R36").Value2 = 6763210
$ws.Range("Z36").Value2 = "11:50"
$ws.Range("AB36").Value2 = "11:50"

# --- Row 37 ---
$ws.Range("A37").Value2 = 130789508
$ws.Range("Q37").Value2 = 490700
$ws.Range("R37").Value2 = 6763430
$ws.Range("Z37").Value2 = "12:07"
$ws.Range("AB37").Value2 = "12:07"

# --- Row 43 ---
$ws.Range("A43").Value2 = 130815695
$ws.Range("B43").Value2 = 57884
$ws.Range("E43").Value2 = 100109
$ws.Range("F43").Value2 = "Tretåig hackspett"
$ws.Range("G43").Value2 = "Picoides tridactylus"
$ws.Range("H43").Value2 = "(Linnaeus, 1758)"
$ws.Range("M43").Value2 = "äldre spår"
$ws.Range("P43").Value2 = "Kråkbackarna, Dlr"
$ws.Range("Q43").Value2 = 490815
$ws.Range("R43").Value2 = 6763446
$ws.Range("S43").Value2 = 10
$ws.Range("AC43").Value2 = "3 bilder på tall"

# --- Row 44 ---
$ws.Range("A44").Value2 = 130807544
$ws.Range("B44").Value2 = 79243
$ws.Range("E44").Value2 = 6425
$ws.Range("F44").Value2 = "Garnlav"
$ws.Range("G44").Value2 = "Alectoria sarmentosa"
$ws.Range("H44").Value2 = "(Ach.) Ach."
$ws.Range("M44").ClearContents() | Out-Null
$ws.Range("P44").Value2 = "Truppan, Dlr"
$ws.Range("Q44").Value2 = 491106
$ws.Range("R44").Value2 = 6763223
$ws.Range("S44").Value2 = 50
$ws.Range("AC44").Value2 = "2 bilder på gran vid basväg samt tall"
